# ============================================================
# Add two METAR lookup sheets (SkyCondition, WeatherType) next
# to the renamed Col_Meteo_US sheet, matching the commit that
# appended "Data/description_colonnes_meteo_US_nov2013_etendu.xlsx".
# ============================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Col_Meteo_US"

# --- New sheet: SkyCondition, right after Col_Meteo_US ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "SkyCondition"

# --- New sheet: WeatherType, right after SkyCondition ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "WeatherType"

# ------------------------------------------------------------
# SkyCondition data (header + 9 METAR sky-cover codes)
# ------------------------------------------------------------
$ws2.Range("A1").Value = "Code"
$ws2.Range("B1").Value = "Signification"
$ws2.Range("C1").Value = "Description"
$ws2.Cells.Item(2, 1).Value = "SKC"
$ws2.Cells.Item(2, 2).Value = "Sky Clear"
$ws2.Cells.Item(2, 3).Value = "Pas de nuages observés"
$ws2.Cells.Item(3, 1).Value = "CLR"
$ws2.Cells.Item(3, 2).Value = "Clear"
$ws2.Cells.Item(3, 3).Value = "Ciel clair (aucun nuage en dessous de 12 000 pieds, station auto)"
$ws2.Cells.Item(4, 1).Value = "FEWnnn"
$ws2.Cells.Item(4, 2).Value = "Few clouds"
$ws2.Cells.Item(4, 3).Value = "1 à 2 octas de couverture (nuages rares), nnn = altitude en centaines de pieds"
$ws2.Cells.Item(5, 1).Value = "SCTnnn"
$ws2.Cells.Item(5, 2).Value = "Scattered clouds"
$ws2.Cells.Item(5, 3).Value = "3 à 4 octas de couverture (nuages épars)"
$ws2.Cells.Item(6, 1).Value = "BKNnnn"
$ws2.Cells.Item(6, 2).Value = "Broken clouds"
$ws2.Cells.Item(6, 3).Value = "5 à 7 octas (nuages fragmentés - couverture significative)"
$ws2.Cells.Item(7, 1).Value = "OVCnnn"
$ws2.Cells.Item(7, 2).Value = "Overcast"
$ws2.Cells.Item(7, 3).Value = "8 octas (ciel couvert)"
$ws2.Cells.Item(8, 1).Value = "VVnnn"
$ws2.Cells.Item(8, 2).Value = "Vertical visibility"
$ws2.Cells.Item(8, 3).Value = "Visibilité verticale limitée (ex. brouillard épais)"
$ws2.Cells.Item(9, 1).Value = "NSC"
$ws2.Cells.Item(9, 2).Value = "No significant clouds"
$ws2.Cells.Item(9, 3).Value = "Aucun nuage significatif"
$ws2.Cells.Item(10, 1).Value = "///"
$ws2.Cells.Item(10, 2).Value = "Indisponible"
$ws2.Cells.Item(10, 3).Value = "Données manquantes ou non valides"
# Copy the header style (bold, centered, bordered) from Col_Meteo_US's
# own header row so the new table matches it exactly, without creating
# any new cell-format entries.
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------
# WeatherType data (header + 26 METAR weather-phenomenon codes)
# ------------------------------------------------------------
$ws3.Range("A1").Value = "Code"
$ws3.Range("B1").Value = "Signification"
$ws3.Range("C1").Value = "Description"
$ws3.Cells.Item(2, 1).Value = "-"
$ws3.Cells.Item(2, 2).Value = "Faible"
$ws3.Cells.Item(2, 3).Value = "Intensité faible du phénomène"
$ws3.Cells.Item(3, 1).Value = "+"
$ws3.Cells.Item(3, 2).Value = "Fort"
$ws3.Cells.Item(3, 3).Value = "Intensité forte"
$ws3.Cells.Item(4, 1).Value = "VC"
$ws3.Cells.Item(4, 2).Value = "Dans les environs"
$ws3.Cells.Item(4, 3).Value = "Phénomène proche (dans un rayon de ~8 km)"
$ws3.Cells.Item(5, 1).Value = "RA"
$ws3.Cells.Item(5, 2).Value = "Rain"
$ws3.Cells.Item(5, 3).Value = "Pluie"
$ws3.Cells.Item(6, 1).Value = "SN"
$ws3.Cells.Item(6, 2).Value = "Snow"
$ws3.Cells.Item(6, 3).Value = "Neige"
$ws3.Cells.Item(7, 1).Value = "SG"
$ws3.Cells.Item(7, 2).Value = "Snow Grains"
$ws3.Cells.Item(7, 3).Value = "Grains de neige"
$ws3.Cells.Item(8, 1).Value = "IC"
$ws3.Cells.Item(8, 2).Value = "Ice Crystals"
$ws3.Cells.Item(8, 3).Value = "Cristaux de glace"
$ws3.Cells.Item(9, 1).Value = "PL"
$ws3.Cells.Item(9, 2).Value = "Ice Pellets"
$ws3.Cells.Item(9, 3).Value = "Granules de glace"
$ws3.Cells.Item(10, 1).Value = "GR"
$ws3.Cells.Item(10, 2).Value = "Hail"
$ws3.Cells.Item(10, 3).Value = "Grêle (≥ 5 mm)"
$ws3.Cells.Item(11, 1).Value = "GS"
$ws3.Cells.Item(11, 2).Value = "Small Hail"
$ws3.Cells.Item(11, 3).Value = "Petite grêle ou granules de neige"
$ws3.Cells.Item(12, 1).Value = "UP"
$ws3.Cells.Item(12, 2).Value = "Unknown Precipitation"
$ws3.Cells.Item(12, 3).Value = "Précipitation non identifiée (station auto)"
$ws3.Cells.Item(13, 1).Value = "DZ"
$ws3.Cells.Item(13, 2).Value = "Drizzle"
$ws3.Cells.Item(13, 3).Value = "Bruine"
$ws3.Cells.Item(14, 1).Value = "FZ"
$ws3.Cells.Item(14, 2).Value = "Freezing"
$ws3.Cells.Item(14, 3).Value = "Préfixe indiquant givre ou verglas"
$ws3.Cells.Item(15, 1).Value = "BR"
$ws3.Cells.Item(15, 2).Value = "Mist"
$ws3.Cells.Item(15, 3).Value = "Brume"
$ws3.Cells.Item(16, 1).Value = "FG"
$ws3.Cells.Item(16, 2).Value = "Fog"
$ws3.Cells.Item(16, 3).Value = "Brouillard"
$ws3.Cells.Item(17, 1).Value = "FU"
$ws3.Cells.Item(17, 2).Value = "Smoke"
$ws3.Cells.Item(17, 3).Value = "Fumée"
$ws3.Cells.Item(18, 1).Value = "VA"
$ws3.Cells.Item(18, 2).Value = "Volcanic Ash"
$ws3.Cells.Item(18, 3).Value = "Cendres volcaniques"
$ws3.Cells.Item(19, 1).Value = "DU"
$ws3.Cells.Item(19, 2).Value = "Dust"
$ws3.Cells.Item(19, 3).Value = "Poussière généralisée"
$ws3.Cells.Item(20, 1).Value = "SA"
$ws3.Cells.Item(20, 2).Value = "Sand"
$ws3.Cells.Item(20, 3).Value = "Sable"
$ws3.Cells.Item(21, 1).Value = "HZ"
$ws3.Cells.Item(21, 2).Value = "Haze"
$ws3.Cells.Item(21, 3).Value = "Brume sèche ou poussiéreuse"
$ws3.Cells.Item(22, 1).Value = "TS"
$ws3.Cells.Item(22, 2).Value = "Thunderstorm"
$ws3.Cells.Item(22, 3).Value = "Orage"
$ws3.Cells.Item(23, 1).Value = "SQ"
$ws3.Cells.Item(23, 2).Value = "Squall"
$ws3.Cells.Item(23, 3).Value = "Rafales violentes"
$ws3.Cells.Item(24, 1).Value = "SS"
$ws3.Cells.Item(24, 2).Value = "Sandstorm"
$ws3.Cells.Item(24, 3).Value = "Tempête de sable"
$ws3.Cells.Item(25, 1).Value = "DS"
$ws3.Cells.Item(25, 2).Value = "Duststorm"
$ws3.Cells.Item(25, 3).Value = "Tempête de poussière"
$ws3.Cells.Item(26, 1).Value = "PO"
$ws3.Cells.Item(26, 2).Value = "Dust/Sand Whirls"
$ws3.Cells.Item(26, 3).Value = "Tourbillons de poussière ou sable"
$ws3.Cells.Item(27, 1).Value = "FC"
$ws3.Cells.Item(27, 2).Value = "Funnel Cloud/Tornado"
$ws3.Cells.Item(27, 3).Value = "Tornade ou trombe"
$ws1.Range("A1:C1").Copy()
$ws3.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------
# Restore Col_Meteo_US as the active sheet, scrolled/selected
# the way it was left in the saved file.
# ------------------------------------------------------------
$ws1.Activate()
$ws1.Range("L37").Select()
